# Add run mode option: new "test_suite" worksheet listing each test case
# and whether it should run.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet (OpenAccountTest).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "test_suite"

# Header row
$ws.Range("A1").Value = "TCID"
$ws.Range("B1").Value = "Runmode"

# Test case rows
$ws.Range("A2").Value = "BankManagerLoginTest"
$ws.Range("B2").Value = "Y"

$ws.Range("A3").Value = "AddCustomerTest"
$ws.Range("B3").Value = "Y"

$ws.Range("A4").Value = "OpenAccountTest"
$ws.Range("B4").Value = "N"

# Give the test case names a distinct (black) font, matching the style
# applied to the TCID column entries.
$ws.Range("A2:A4").Font.Color = 0

# Make the new sheet the active tab/selection, same as when it was added
# and last edited in Excel.
[void]$ws.Activate()
$ws.Range("B4").Select() | Out-Null
